$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 112223201
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "2"
$ws.Range("I4").Style = "Normal"
$ws.Range("Q4").Value = 636549
$ws.Range("R4").Value = 6523814
$ws.Range("Z4").Value = "21:15"
$ws.Range("AB4").Value = "21:15"
$ws.Range("AC4").Value = "Två hanar observerade med ljud och visuellt."

# Row 5
$ws.Range("A5").Value = 112223193
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "5"
$ws.Range("I5").Style = "Normal"
$ws.Range("Q5").Value = 636408
$ws.Range("R5").Value = 6524025
$ws.Range("AC5").Value = "Hanar i lekdräkt."

# Row 6
$ws.Range("A6").Value = 112223203
$ws.Range("J6").Value = "ex."
$ws.Range("K6").Value = "adult"
$ws.Range("L6").Value = "hane"
$ws.Range("M6").Value = "spel/sång"
$ws.Range("Q6").Value = 636549
$ws.Range("R6").Value = 6523814
$ws.Range("Z6").Value = "21:15"
$ws.Range("AB6").Value = "21:15"
$ws.Range("AC6").Value = "Observerad med ljud."

# Row 7
$ws.Range("A7").Value = 112223199
$ws.Range("B7").Value = 57584
$ws.Range("E7").Value = 208245
$ws.Range("F7").Value = "Vanlig padda"
$ws.Range("G7").Value = "Bufo bufo"
$ws.Range("M7").Value = "spel/sång"
$ws.Range("Q7").Value = 636422
$ws.Range("R7").Value = 6523909
$ws.Range("Z7").Value = "20:30"
$ws.Range("AB7").Value = "20:30"
$ws.Range("AC7").Value = "Observerad med ljud."

# Row 8
$ws.Range("A8").Value = 112223196
$ws.Range("B8").Value = 57620
$ws.Range("E8").Value = 208242
$ws.Range("F8").Value = "Mindre vattensalamander"
$ws.Range("G8").Value = "Lissotriton vulgaris"
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = "1"
$ws.Range("I8").Style = "Normal"
$ws.Range("M8").Value = "i vatten/simmande"
$ws.Range("Q8").Value = 636422
$ws.Range("R8").Value = 6523909
$ws.Range("Z8").Value = "20:30"
$ws.Range("AB8").Value = "20:30"
$ws.Range("AC8").Value = "Hane i lekdräkt."

# Row 9
$ws.Range("A9").Value = 112223188
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "1"
$ws.Range("I9").Style = "Normal"
$ws.Range("Q9").Value = 636399
$ws.Range("R9").Value = 6523963
$ws.Range("Z9").Value = "20:20"
$ws.Range("AB9").Value = "20:20"
$ws.Range("AC9").ClearContents()

# Row 10
$ws.Range("A10").Value = 112223184
$ws.Range("J10").Value = "äggklumpar"
$ws.Range("K10").Value = "ägg"
$ws.Range("L10").Value = ""
$ws.Range("M10").Value = ""
$ws.Range("Q10").Value = 636399
$ws.Range("R10").Value = 6523963
$ws.Range("Z10").Value = "20:20"
$ws.Range("AB10").Value = "20:20"
$ws.Range("AC10").Value = "Romklump. Troligen åkergroda som förekommer i närliggande dammar."
